# Update handoff/handback datetimes on the zh-cn and de-de report sheets
# as part of regenerating the handback status report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-08 11:49:20"
$wsZhCn.Range("G2").Value = "2016-01-08 11:50:08"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-08 11:49:33"
$wsDeDe.Range("G2").Value = "2016-01-08 11:50:29"
